$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a range to be treated as text so that numeric-looking
# strings (e.g. "1.00", "0.999") are not coerced into numbers by Excel,
# matching the original inlineStr cell type used throughout the sheet.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "95.111.18"
$ws.Range("E2").Value = "  -1.49%  "

Set-TextValue $ws.Range("D3") "3.445.63"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "239.64"
$ws.Range("E5").Value = "  -4.10%  "

Set-TextValue $ws.Range("D6") "642.29"
$ws.Range("E6").Value = "  -1.94%  "

Set-TextValue $ws.Range("D7") "1.47"
$ws.Range("E7").Value = "  +5.18%  "

Set-TextValue $ws.Range("D8") "0.402"
$ws.Range("E8").Value = "  -4.42%  "

$ws.Range("E9").Value = "  +0.04%  "

Set-TextValue $ws.Range("D10") "0.994"
$ws.Range("E10").Value = "  +0.25%  "

Set-TextValue $ws.Range("D11") "3.442.87"
$ws.Range("E11").Value = "  +3.51%  "

$ws.Range("E12").Value = "  -4.13%  "

Set-TextValue $ws.Range("D13") "41.49"
$ws.Range("E13").Value = "  +3.06%  "

Set-TextValue $ws.Range("D14") "6.08"
$ws.Range("E14").Value = "  +0.02%  "

Set-TextValue $ws.Range("D15") "94.895.10"
$ws.Range("E15").Value = "  -1.43%  "

Set-TextValue $ws.Range("D16") "4.092.56"
$ws.Range("E16").Value = "  +3.46%  "

$ws.Range("E17").Value = "  +2.00%  "

Set-TextValue $ws.Range("D18") "8.43"
$ws.Range("E18").Value = "  -1.03%  "

Set-TextValue $ws.Range("D19") "3.437.89"
$ws.Range("E19").Value = "  +3.29%  "

Set-TextValue $ws.Range("D20") "17.70"
$ws.Range("E20").Value = "  +3.60%  "

Set-TextValue $ws.Range("D21") "11.43"
$ws.Range("E21").Value = "  +8.58%  "

Set-TextValue $ws.Range("D22") "0.505"
$ws.Range("E22").Value = "  -6.91%  "

Set-TextValue $ws.Range("D23") "501.69"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("E24").Value = "  -5.92%  "

Set-TextValue $ws.Range("D25") "0.0000191"
$ws.Range("E25").Value = "  -3.05%  "

Set-TextValue $ws.Range("D26") "6.55"
$ws.Range("E26").Value = "  -0.40%  "

Set-TextValue $ws.Range("D27") "91.82"
$ws.Range("E27").Value = "  -4.60%  "

Set-TextValue $ws.Range("D28") "3.628.78"
$ws.Range("E28").Value = "  +3.43%  "

Set-TextValue $ws.Range("D29") "11.98"
$ws.Range("E29").Value = "  -0.84%  "

Set-TextValue $ws.Range("D30") "11.68"
$ws.Range("E30").Value = "  +5.71%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("E33").Value = "  -4.78%  "

$ws.Range("E34").Value = "  -1.70%  "

Set-TextValue $ws.Range("D35") "31.01"
$ws.Range("E35").Value = "  +11.00%  "

Set-TextValue $ws.Range("D36") "0.999"
$ws.Range("E36").Value = "  -0.13%  "

Set-TextValue $ws.Range("D37") "0.565"
$ws.Range("E37").Value = "  +3.68%  "

Set-TextValue $ws.Range("D38") "7.66"
$ws.Range("E38").Value = "  +0.62%  "

Set-TextValue $ws.Range("D39") "1.44"
$ws.Range("E39").Value = "  -2.26%  "

Set-TextValue $ws.Range("D40") "522.02"
$ws.Range("E40").Value = "  +3.05%  "

$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("E42").Value = "  -0.54%  "

Set-TextValue $ws.Range("D43") "0.910"
$ws.Range("E43").Value = "  +9.70%  "

Set-TextValue $ws.Range("D44") "24.08"
$ws.Range("E44").Value = "  -1.19%  "

$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0415"
$ws.Range("E46").Value = "  -2.30%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D47") "5.60"
$ws.Range("E47").Value = "  +1.70%  "

Set-TextValue $ws.Range("D48") "3.49"
$ws.Range("E48").Value = "  -4.61%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D49") "2.14"
$ws.Range("E49").Value = "  +8.64%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D50") "53.57"
$ws.Range("E50").Value = "  +0.46%  "

Set-TextValue $ws.Range("D51") "3.18"
$ws.Range("E51").Value = "  +1.99%  "
